$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.45085973298193
$ws.Range("C2").Value = 2.87357769895942

$ws.Range("B3").Value = 1.0766438326074
$ws.Range("C3").Value = 0.767058349939147
$ws.Range("E3").Value = 0.966

$ws.Range("B4").Value = 2.91352809337216
$ws.Range("C4").Value = 3.46168547710751

$ws.Range("B5").Value = 1.32184967903658
$ws.Range("C5").Value = 1.03432976319099
$ws.Range("E5").Value = 1.152

$ws.Range("B6").Value = 1.18955792454387
$ws.Range("C6").Value = 0.859637113798726
$ws.Range("E6").Value = 1.057

$ws.Range("B20").Select()
